# Updates the cryptocurrency price/volume snapshot values in the
# active worksheet to reflect the latest scrape (GitHub Actions run).
# Also corrects the ranking order of two coin pairs whose positions
# were swapped (Uniswap/Avalanche and TrustWalletToken/HuobiToken).
#
# Some "Price" cells hold numeric-looking text (e.g. "212.36") that must
# stay plain text, matching the source data's inline-string cell type.
# A leading apostrophe forces Excel to keep such values as text instead
# of auto-converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.908.00'
$ws.Range("E2").Value = '  -0.09%  '
$ws.Range("D3").Value = '1.636.15'
$ws.Range("E3").Value = '  -0.71%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''212.36'
$ws.Range("E5").Value = '  -0.58%  '
$ws.Range("D6").Value = '''0.523'
$ws.Range("E6").Value = '  -0.27%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '''23.27'
$ws.Range("E8").Value = '  -1.23%  '
$ws.Range("E9").Value = '  -2.52%  '
$ws.Range("E10").Value = '  -0.19%  '
$ws.Range("D11").Value = '''0.0880'
$ws.Range("E11").Value = '  +0.98%  '
$ws.Range("D12").Value = '1.869.16'
$ws.Range("E12").Value = '  -0.69%  '
$ws.Range("D13").Value = '1.636.04'
$ws.Range("E13").Value = '  -0.15%  '
$ws.Range("E14").Value = '  -0.42%  '
$ws.Range("D15").Value = '''0.569'
$ws.Range("E15").Value = '  +0.99%  '
$ws.Range("D16").Value = '''65.17'
$ws.Range("E16").Value = '  -0.81%  '
$ws.Range("D17").Value = '27.916.70'
$ws.Range("E17").Value = '  -0.05%  '
$ws.Range("D18").Value = '''230.20'
$ws.Range("E18").Value = '  -0.77%  '
$ws.Range("D19").Value = '0.0₃0720'
$ws.Range("E19").Value = '  -0.38%  '
$ws.Range("D20").Value = '''7.51'
$ws.Range("E20").Value = '  -1.97%  '
$ws.Range("E21").Value = '  -0.07%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '''4.37'
$ws.Range("E22").Value = '  -0.53%  '
$ws.Range("B23").Value = 'Avalanche'
$ws.Range("C23").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D23").Value = '''10.36'
$ws.Range("E23").Value = '  -3.22%  '
$ws.Range("D24").Value = '''2.07'
$ws.Range("E24").Value = '  -3.57%  '
$ws.Range("D25").Value = '''153.66'
$ws.Range("E25").Value = '  +1.14%  '
$ws.Range("E26").Value = '  +0.81%  '
$ws.Range("E27").Value = '  -0.48%  '
$ws.Range("E28").Value = '  -0.57%  '
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("E30").Value = '  -0.60%  '
$ws.Range("E31").Value = '  -0.77%  '
$ws.Range("D32").Value = '''3.37'
$ws.Range("E32").Value = '  +0.76%  '
$ws.Range("D33").Value = '1.407.90'
$ws.Range("E33").Value = '  -3.19%  '
$ws.Range("E34").Value = '  -1.48%  '
$ws.Range("E35").Value = '  +1.22%  '
$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D36").Value = '''1.01'
$ws.Range("E36").Value = '  +10.18%  '
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").Value = '''2.36'
$ws.Range("E37").Value = '  +1.56%  '
$ws.Range("E38").Value = '  +0.45%  '
$ws.Range("E39").Value = '  -0.14%  '
$ws.Range("D40").Value = '''0.874'
$ws.Range("E40").Value = '  -1.72%  '
$ws.Range("E41").Value = '  +0.22%  '
$ws.Range("D43").Value = '''67.04'
$ws.Range("E43").Value = '  -3.39%  '
$ws.Range("D44").Value = '''5.53'
$ws.Range("E44").Value = '  +2.60%  '
$ws.Range("D45").Value = '''1.81'
$ws.Range("E45").Value = '  +1.44%  '
$ws.Range("E46").Value = '  -1.75%  '
$ws.Range("D47").Value = '1.778.43'
$ws.Range("E47").Value = '  -0.75%  '
$ws.Range("D48").Value = '''87.93'
$ws.Range("E48").Value = '  -1.20%  '
$ws.Range("E49").Value = '  -0.36%  '
$ws.Range("E50").Value = '  -0.26%  '
$ws.Range("D51").Value = '''7.58'
$ws.Range("E51").Value = '  -1.69%  '
